$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.363.61"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.590.67"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "'210.26"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.0610"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "1.593.82"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "'64.38"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "26.372.92"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("D20").Value = "'210.80"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "'8.92"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'145.00"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "'15.25"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "'3.21"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "1.308.41"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "'0.614"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  -13.20%  "
$ws.Range("D40").Value = "'0.810"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D42").Value = "'5.63"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "1.728.22"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("D50").Value = "'0.0504"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  -0.46%  "
